# semana 45 de 2025
# Adds week 44 (column AU) and week 45 (column AV) data to the weekly
# IRA-hospital revision sheet, plus a couple of small corrections that
# came in with this week's data refresh (row 31 counts, row 52 UPGD name).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row: new week-number columns AU1 ("44") and AV1 ("45").
# These must land as TEXT (matching the existing "1".."43" header
# cells), so we use the quote-prefix trick to force text entry, then
# copy the number format / font / alignment from AT1 (the existing
# "43" header) onto the two new header cells so they end up styled
# identically to the rest of the header row.
# ---------------------------------------------------------------------
$ws.Range("AU1").Value = "'44"
$ws.Range("AV1").Value = "'45"
$ws.Range("AT1").Copy()
$ws.Range("AU1:AV1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# New weekly-count data for columns AU (week 44) and AV (week 45).
# Only rows that already carried data through column AT receive the
# new cells, matching the source report's per-row sparsity.
# ---------------------------------------------------------------------
$ws.Range("AU2").Value = 0
$ws.Range("AV2").Value = 0
$ws.Range("AV3").Value = 0
$ws.Range("AU5").Value = 0
$ws.Range("AV5").Value = 0
$ws.Range("AU6").Value = 18
$ws.Range("AV6").Value = 33
$ws.Range("AU7").Value = 18
$ws.Range("AV7").Value = 0
$ws.Range("AU8").Value = 5
$ws.Range("AV8").Value = 8
$ws.Range("AU9").Value = 0
$ws.Range("AV9").Value = 0
$ws.Range("AU10").Value = 0
$ws.Range("AU12").Value = 0
$ws.Range("AV12").Value = 0
$ws.Range("AV13").Value = 0
$ws.Range("AU14").Value = 0
$ws.Range("AV14").Value = 0
$ws.Range("AU15").Value = 0
$ws.Range("AV15").Value = 0
$ws.Range("AV16").Value = 0
$ws.Range("AV17").Value = 0
$ws.Range("AV22").Value = 0
$ws.Range("AU23").Value = 0
$ws.Range("AV23").Value = 0
$ws.Range("AU24").Value = 0
$ws.Range("AU25").Value = 1
$ws.Range("AV25").Value = 4
$ws.Range("AV26").Value = 0
$ws.Range("AU28").Value = 4
$ws.Range("AV28").Value = 5
$ws.Range("AU29").Value = 3
$ws.Range("AV29").Value = 2
$ws.Range("AU30").Value = 2
$ws.Range("AV30").Value = 3
$ws.Range("AU31").Value = 0
$ws.Range("AV31").Value = 0
$ws.Range("AU34").Value = 0
$ws.Range("AU35").Value = 2
$ws.Range("AV35").Value = 5
$ws.Range("AU36").Value = 0
$ws.Range("AV36").Value = 0
$ws.Range("AU37").Value = 0
$ws.Range("AV37").Value = 0
$ws.Range("AU38").Value = 0
$ws.Range("AV38").Value = 0
$ws.Range("AU40").Value = 0
$ws.Range("AT41").Value = 0
$ws.Range("AU41").Value = 0
$ws.Range("AV41").Value = 0
$ws.Range("AU42").Value = 0
$ws.Range("AV42").Value = 0
$ws.Range("AU43").Value = 0
$ws.Range("AV43").Value = 0
$ws.Range("AU44").Value = 0
$ws.Range("AU45").Value = 0
$ws.Range("AV45").Value = 0
$ws.Range("AU46").Value = 0
$ws.Range("AV46").Value = 0
$ws.Range("AU47").Value = 0
$ws.Range("AV47").Value = 0
$ws.Range("AU48").Value = 0
$ws.Range("AV48").Value = 0
$ws.Range("AU49").Value = 0
$ws.Range("AV49").Value = 0
$ws.Range("AU50").Value = 0
$ws.Range("AV50").Value = 0
$ws.Range("AU51").Value = 0
$ws.Range("AU53").Value = 0
$ws.Range("AV53").Value = 0
$ws.Range("AU54").Value = 0
$ws.Range("AV54").Value = 0
$ws.Range("AU55").Value = 0
$ws.Range("AV55").Value = 0
$ws.Range("AU56").Value = 0
$ws.Range("AV56").Value = 0
$ws.Range("AU57").Value = 0
$ws.Range("AV57").Value = 0
$ws.Range("AU58").Value = 0
$ws.Range("AV58").Value = 0

# ---------------------------------------------------------------------
# Row 31 (FUNDACIÓN LA LIGA AMA SALVAR VIDAS) - revised weekly counts.
# ---------------------------------------------------------------------
$ws.Range("N31").Value = 1
$ws.Range("T31").Value = 1
$ws.Range("AJ31").Value = 1
$ws.Range("AM31").Value = 1

# ---------------------------------------------------------------------
# Row 52 - UPGD name correction.
# ---------------------------------------------------------------------
$ws.Range("C52").Value = "COOMEVA EXPERIENCIA MEDICA SAS"
